# The "sentiments" column (D) was removed entirely from the sheet: the
# header in D1 and every sentiment value in D2:D89 are deleted, and the
# remaining columns (A:C) shift to fill the sheet dimension (A1:C89).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").EntireColumn.Delete()
